$d = $word.ActiveDocument
$r = $d.Range(0, 0)

$frag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:spacing w:before="180" w:after="180"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
<w:color w:val="000000"/>
</w:rPr>
<w:t xml:space="preserve">This dataset describes the palatability, for the land crab </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
<w:i/>
<w:iCs/>
<w:color w:val="000000"/>
</w:rPr>
<w:t>Cardisoma</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
<w:i/>
<w:iCs/>
<w:color w:val="000000"/>
</w:rPr>
<w:t xml:space="preserve"> carnifex</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
<w:color w:val="000000"/>
</w:rPr>
<w:t xml:space="preserve">, of two rodenticide bait products available for conservation use in the U.S.: “25W” containing 25 ppm brodifacoum as the active ingredient, and “D50” containing 50 ppm </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
<w:color w:val="000000"/>
</w:rPr>
<w:t>diphacinone</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Cambria" w:eastAsia="Times New Roman" w:hAnsi="Cambria" w:cs="Times New Roman"/>
<w:color w:val="000000"/>
</w:rPr>
<w:t xml:space="preserve"> as the active ingredient. Both rodenticide bait products were found to be palatable to rats and crabs when presented alongside three commonly available food items: coconut endosperm (meat), the meristematic tissue of young coconut palms, and the fleshy mesocarp of Pandanus fruit. </w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$r.InsertXML($frag)
